$d = $word.ActiveDocument

# Locate the paragraph to replace: "SQL> insert into department values(dept_id_seq.nextval,'CSE',51004);"
# (the final content paragraph, immediately preceded/followed by an empty paragraph). Track its
# 1-based index manually -- Paragraph.Index is not re-evaluated once the document is mutated.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "SQL> insert into department values(dept_id_seq.nextval,'CSE',51004);*") {
        $targetIndex = $i
    }
}

$target = $d.Paragraphs.Item($targetIndex)
$r = $target.Range

# Add a new blank paragraph, then a "Q2)" paragraph, right before the target paragraph.
$r.InsertParagraphBefore()
$r.InsertParagraphBefore()

# The target paragraph itself is now two slots further down; the blank paragraph immediately
# before it is the freshly inserted "Q2)" placeholder.
$targetIndex = $targetIndex + 2
$qPara = $d.Paragraphs.Item($targetIndex - 1)
$qPara.Range.Text = "Q2)"

# Replace the target paragraph's content with the "drop sequence" statement, preserving the
# spell-check proofErr markup around the identifier, exactly like the rest of the document.
$target = $d.Paragraphs.Item($targetIndex)
$dropXml = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">SQL&gt; drop sequence </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>dept_id_seq</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>;</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$target.Range.InsertXML($dropXml)

# After it: a blank paragraph, then "Sequence dropped." (the trailing blank paragraph that used
# to follow the old "insert into department" paragraph is left in place, after these new ones).
$target = $d.Paragraphs.Item($targetIndex)
$r2 = $target.Range
$r2.InsertParagraphAfter()
$r2.InsertParagraphAfter()

$droppedPara = $d.Paragraphs.Item($targetIndex + 2)
$droppedPara.Range.Text = "Sequence dropped."
